# Auto-generated script to append scrim result rows (Equipo 1/2 match history)
$wb = $excel.ActiveWorkbook

# --- Worksheet index 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("A4:N4").Copy($ws.Range("A43:N43"))
$vals = @('CHARLIE', 'CORDELIUS', 'BARLEY', 'MEG', 'JUJU', 'CROW', 'Equipo 1', 'FUT|Nowy297', 'FUT|MeOw', 'FUT|GeRo', 'TH|LeNain', 'TH|iKaoss', 'TH|Zhar', '20250724T165035.000Z')
for ($i = 0; $i -lt 14; $i++) { $ws.Cells.Item(43, $i+1).Value = $vals[$i] }

$ws.Range("A4:N4").Copy($ws.Range("A44:N44"))
$vals = @('CHARLIE', 'CORDELIUS', 'BARLEY', 'MEG', 'JUJU', 'CROW', 'Equipo 1', 'FUT|Nowy297', 'FUT|MeOw', 'FUT|GeRo', 'TH|LeNain', 'TH|iKaoss', 'TH|Zhar', '20250724T164730.000Z')
for ($i = 0; $i -lt 14; $i++) { $ws.Cells.Item(44, $i+1).Value = $vals[$i] }

$ws.Range("A4:N4").Copy($ws.Range("A45:N45"))
$vals = @('CHARLIE', 'MEG', 'JUJU', 'KAZE', 'CORDELIUS', 'BEA', 'Equipo 1', 'FUT|GeRo', 'FUT|Nowy297', 'FUT|MeOw', 'TH|LeNain', 'TH|iKaoss', 'TH|Zhar', '20250724T164056.000Z')
for ($i = 0; $i -lt 14; $i++) { $ws.Cells.Item(45, $i+1).Value = $vals[$i] }

$ws.Range("A5:N5").Copy($ws.Range("A46:N46"))
$vals = @('CHARLIE', 'MEG', 'JUJU', 'KAZE', 'CORDELIUS', 'BEA', 'Equipo 2', 'FUT|GeRo', 'FUT|Nowy297', 'FUT|MeOw', 'TH|LeNain', 'TH|iKaoss', 'TH|Zhar', '20250724T163801.000Z')
for ($i = 0; $i -lt 14; $i++) { $ws.Cells.Item(46, $i+1).Value = $vals[$i] }

$ws.Range("A4:N4").Copy($ws.Range("A47:N47"))
$vals = @('CHARLIE', 'MEG', 'JUJU', 'KAZE', 'CORDELIUS', 'BEA', 'Equipo 1', 'FUT|GeRo', 'FUT|Nowy297', 'FUT|MeOw', 'TH|LeNain', 'TH|iKaoss', 'TH|Zhar', '20250724T163442.000Z')
for ($i = 0; $i -lt 14; $i++) { $ws.Cells.Item(47, $i+1).Value = $vals[$i] }


# --- Worksheet index 5 ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("A4:N4").Copy($ws.Range("A60:N60"))
$vals = @('GUS', 'BUZZ', 'SQUEAK', 'SHADE', 'MEEPLE', 'CARL', 'Equipo 2', 'SK|Ope', 'SK|Joker', 'SK|Yoshi825', 'HMB|Lukii', 'HMB|Symantec', 'HMB|BosS', '20250724T164723.000Z')
for ($i = 0; $i -lt 14; $i++) { $ws.Cells.Item(60, $i+1).Value = $vals[$i] }

$ws.Range("A4:N4").Copy($ws.Range("A61:N61"))
$vals = @('GUS', 'BUZZ', 'SQUEAK', 'SHADE', 'MEEPLE', 'CARL', 'Equipo 2', 'SK|Ope', 'SK|Joker', 'SK|Yoshi825', 'HMB|Lukii', 'HMB|Symantec', 'HMB|BosS', '20250724T164518.000Z')
for ($i = 0; $i -lt 14; $i++) { $ws.Cells.Item(61, $i+1).Value = $vals[$i] }

$ws.Range("A4:N4").Copy($ws.Range("A62:N62"))
$vals = @('GRAY', 'JUJU', 'ANGELO', 'SQUEAK', 'MEEPLE', 'WILLOW', 'Equipo 2', 'SK|Ope', 'SK|Yoshi825', 'SK|Joker', 'HMB|BosS', 'HMB|Lukii', 'HMB|Symantec', '20250724T163927.000Z')
for ($i = 0; $i -lt 14; $i++) { $ws.Cells.Item(62, $i+1).Value = $vals[$i] }

$ws.Range("A4:N4").Copy($ws.Range("A63:N63"))
$vals = @('GRAY', 'JUJU', 'ANGELO', 'SQUEAK', 'MEEPLE', 'WILLOW', 'Equipo 2', 'SK|Ope', 'SK|Yoshi825', 'SK|Joker', 'HMB|BosS', 'HMB|Lukii', 'HMB|Symantec', '20250724T163707.000Z')
for ($i = 0; $i -lt 14; $i++) { $ws.Cells.Item(63, $i+1).Value = $vals[$i] }

$ws.Range("A9:N9").Copy($ws.Range("A64:N64"))
$vals = @('GRAY', 'JUJU', 'ANGELO', 'SQUEAK', 'MEEPLE', 'WILLOW', 'Equipo 1', 'SK|Ope', 'SK|Yoshi825', 'SK|Joker', 'HMB|BosS', 'HMB|Lukii', 'HMB|Symantec', '20250724T163447.000Z')
for ($i = 0; $i -lt 14; $i++) { $ws.Cells.Item(64, $i+1).Value = $vals[$i] }


# --- Worksheet index 6 ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("A4:N4").Copy($ws.Range("A58:N58"))
$vals = @('BEA', 'CORDELIUS', 'SHADE', 'SURGE', 'LUMI', 'ASH', 'Equipo 2', 'SK|Ope', 'SK|Yoshi825', 'SK|Joker', 'HMB|BosS', 'HMB|Lukii', 'HMB|Symantec', '20250724T162646.000Z')
for ($i = 0; $i -lt 14; $i++) { $ws.Cells.Item(58, $i+1).Value = $vals[$i] }

$ws.Range("A4:N4").Copy($ws.Range("A59:N59"))
$vals = @('BEA', 'CORDELIUS', 'SHADE', 'SURGE', 'LUMI', 'ASH', 'Equipo 2', 'SK|Ope', 'SK|Yoshi825', 'SK|Joker', 'HMB|BosS', 'HMB|Lukii', 'HMB|Symantec', '20250724T162502.000Z')
for ($i = 0; $i -lt 14; $i++) { $ws.Cells.Item(59, $i+1).Value = $vals[$i] }

$ws.Range("A10:N10").Copy($ws.Range("A60:N60"))
$vals = @('BEA', 'CORDELIUS', 'SHADE', 'SURGE', 'LUMI', 'ASH', 'Equipo 1', 'SK|Ope', 'SK|Yoshi825', 'SK|Joker', 'HMB|BosS', 'HMB|Lukii', 'HMB|Symantec', '20250724T162244.000Z')
for ($i = 0; $i -lt 14; $i++) { $ws.Cells.Item(60, $i+1).Value = $vals[$i] }

$ws.Range("A4:N4").Copy($ws.Range("A61:N61"))
$vals = @('LOU', 'HANK', 'KENJI', 'ASH', 'CORDELIUS', 'BEA', 'Equipo 2', 'FUT|GeRo', 'FUT|MeOw', 'FUT|Nowy297', 'TH|LeNain', 'TH|iKaoss', 'TH|Zhar', '20250724T162726.000Z')
for ($i = 0; $i -lt 14; $i++) { $ws.Cells.Item(61, $i+1).Value = $vals[$i] }

$ws.Range("A10:N10").Copy($ws.Range("A62:N62"))
$vals = @('LOU', 'HANK', 'KENJI', 'ASH', 'CORDELIUS', 'BEA', 'Equipo 1', 'FUT|GeRo', 'FUT|MeOw', 'FUT|Nowy297', 'TH|LeNain', 'TH|iKaoss', 'TH|Zhar', '20250724T162536.000Z')
for ($i = 0; $i -lt 14; $i++) { $ws.Cells.Item(62, $i+1).Value = $vals[$i] }

$ws.Range("A4:N4").Copy($ws.Range("A63:N63"))
$vals = @('LOU', 'HANK', 'KENJI', 'ASH', 'CORDELIUS', 'BEA', 'Equipo 2', 'FUT|GeRo', 'FUT|MeOw', 'FUT|Nowy297', 'TH|LeNain', 'TH|iKaoss', 'TH|Zhar', '20250724T162252.000Z')
for ($i = 0; $i -lt 14; $i++) { $ws.Cells.Item(63, $i+1).Value = $vals[$i] }

